# Update "2024" sheet: push a new September log entry onto the stack.
# Every row from 35 downward shifts down by one (the sheet is a
# chronological log, newest first), and the vacated row 35 receives the
# brand-new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Shift rows 35:89 down to 36:90, preserving all columns (A-Y), then
# populate the freed row 35 with the new log entry.
$ws.Rows.Item(35).Insert()

$ws.Range("R35").Value = "bal axis"
$ws.Range("S35").Value = "2024-09-07 08:46:40"
